$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H2").Value = 466.66666
$ws.Range("J2").Value = 666.6667
$ws.Range("L2").Value = 666.6667
$ws.Range("N2").Value = -892.6667
$ws.Range("H32").Value = 1500
$ws.Range("J32").Value = 1500
$ws.Range("L32").Value = 1500
$ws.Range("N32").Value = -2152
$ws.Range("H33").Value = 111.5625
$ws.Range("J33").Value = 97.55556
$ws.Range("L33").Value = 97.55556
$ws.Range("N33").Value = -555.55556
$ws.Range("H51").Value = 4250
$ws.Range("J51").Value = 4250
$ws.Range("L51").Value = 4250
$ws.Range("N51").Value = -5218
$ws.Range("H98").Value = 3991
$ws.Range("I98").Value = 3613.75
$ws.Range("K98").Value = 3613.75
$ws.Range("M98").Value = -2115.75
$ws.Range("H100").Value = 1698.3846
$ws.Range("I100").Value = 1756.5834
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 1756.5834
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -1215.5834
$ws.Range("N100").Value = -2082
$ws.Range("H122").Value = 3991
$ws.Range("I122").Value = 3613.75
$ws.Range("K122").Value = 10841.25
$ws.Range("M122").Value = -8391.25
$ws.Range("H127").Value = 1898.1333
$ws.Range("I127").Value = 1924.4166
$ws.Range("K127").Value = 5773.2498
$ws.Range("M127").Value = -813.2497999999996
$ws.Range("H129").Value = 880.10205
$ws.Range("I129").Value = 774.25
$ws.Range("J129").Value = 889.5111000000001
$ws.Range("K129").Value = 2322.75
$ws.Range("L129").Value = 2668.5333
$ws.Range("M129").Value = 2677.25
$ws.Range("N129").Value = -12668.5333
$ws.Range("H132").Value = 1061.5
$ws.Range("I132").Value = 991.7778
$ws.Range("J132").Value = 1438
$ws.Range("K132").Value = 2975.3334
$ws.Range("L132").Value = 4314
$ws.Range("M132").Value = -445.3334
$ws.Range("N132").Value = -9374
$ws.Range("H137").Value = 2142.8667
$ws.Range("I137").Value = 967.3333
$ws.Range("J137").Value = 2436.75
$ws.Range("K137").Value = 2901.9999
$ws.Range("L137").Value = 7310.25
$ws.Range("M137").Value = -351.9998999999998
$ws.Range("N137").Value = -12410.25
$ws.Range("H138").Value = 2488.2593
$ws.Range("I138").Value = 2921.8572
$ws.Range("J138").Value = 2212.3333
$ws.Range("K138").Value = 8765.571599999999
$ws.Range("L138").Value = 6636.999899999999
$ws.Range("M138").Value = -3625.571599999999
$ws.Range("N138").Value = -16916.9999

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 2225.0762
$ws.Range("I32").Value = 1523.7073
$ws.Range("K32").Value = 1523.7073
$ws.Range("M32").Value = -1236.7073
$ws.Range("H61").Value = 6801.8335
$ws.Range("I61").Value = 2728
$ws.Range("K61").Value = 2728
$ws.Range("M61").Value = -2516
$ws.Range("H74").Value = 1267.7931
$ws.Range("I74").Value = 843.9
$ws.Range("J74").Value = 2209.7778
$ws.Range("K74").Value = 843.9
$ws.Range("L74").Value = 2209.7778
$ws.Range("M74").Value = 30.10000000000002
$ws.Range("N74").Value = -3957.7778
$ws.Range("H77").Value = 1267.7931
$ws.Range("I77").Value = 843.9
$ws.Range("J77").Value = 2209.7778
$ws.Range("K77").Value = 4219.5
$ws.Range("L77").Value = 11048.889
$ws.Range("M77").Value = 148.5
$ws.Range("N77").Value = -19784.889
$ws.Range("H122").Value = 1139.3
$ws.Range("I122").Value = 1318.8462
$ws.Range("K122").Value = 3956.5386
$ws.Range("M122").Value = -1506.5386
$ws.Range("H123").Value = 66665.664
$ws.Range("J123").Value = 66665.664
$ws.Range("L123").Value = 66665.664
$ws.Range("N123").Value = -76465.664
$ws.Range("H132").Value = 1752.0682
$ws.Range("I132").Value = 1176.3823
$ws.Range("J132").Value = 3709.4
$ws.Range("K132").Value = 3529.1469
$ws.Range("L132").Value = 11128.2
$ws.Range("M132").Value = -999.1468999999997
$ws.Range("N132").Value = -16188.2
$ws.Range("H136").Value = 6801.8335
$ws.Range("I136").Value = 2728
$ws.Range("K136").Value = 8184
$ws.Range("M136").Value = -5634

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 1950.0834
$ws.Range("I20").Value = 1871.8334
$ws.Range("K20").Value = 1871.8334
$ws.Range("M20").Value = -1624.8334
$ws.Range("H94").Value = 1010.44446
$ws.Range("J94").Value = 929.6667
$ws.Range("L94").Value = 929.6667
$ws.Range("N94").Value = -1831.6667
$ws.Range("H99").Value = 1886.591
$ws.Range("I99").Value = 1764.6428
$ws.Range("J99").Value = 2100
$ws.Range("K99").Value = 1764.6428
$ws.Range("L99").Value = 2100
$ws.Range("M99").Value = -266.6428000000001
$ws.Range("N99").Value = -5096
$ws.Range("H105").Value = 2253.9333
$ws.Range("I105").Value = 2253.9333
$ws.Range("K105").Value = 2253.9333
$ws.Range("M105").Value = -506.9333000000001
$ws.Range("H108").Value = 94999.5
$ws.Range("J108").Value = 94999.5
$ws.Range("L108").Value = 94999.5
$ws.Range("N108").Value = -102679.5
$ws.Range("H134").Value = 7274.926
$ws.Range("I134").Value = 7674.913
$ws.Range("J134").Value = 4975
$ws.Range("K134").Value = 23024.739
$ws.Range("L134").Value = 14925
$ws.Range("M134").Value = -20489.739
$ws.Range("N134").Value = -19995

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 3069.7222
$ws.Range("I31").Value = 3021.6667
$ws.Range("J31").Value = 3165.8333
$ws.Range("K31").Value = 3021.6667
$ws.Range("L31").Value = 3165.8333
$ws.Range("M31").Value = -2726.6667
$ws.Range("N31").Value = -3755.8333
$ws.Range("H34").Value = 3069.7222
$ws.Range("I34").Value = 3021.6667
$ws.Range("J34").Value = 3165.8333
$ws.Range("K34").Value = 3021.6667
$ws.Range("L34").Value = 3165.8333
$ws.Range("M34").Value = -2819.6667
$ws.Range("N34").Value = -3569.8333
$ws.Range("H99").Value = 2857.8
$ws.Range("I99").Value = 2429.6667
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 2429.6667
$ws.Range("L99").Value = 3500
$ws.Range("M99").Value = -931.6667000000002
$ws.Range("N99").Value = -6496
$ws.Range("H122").Value = 2392.7827
$ws.Range("I122").Value = 938.3333
$ws.Range("K122").Value = 2814.9999
$ws.Range("M122").Value = -364.9998999999998
$ws.Range("H126").Value = 2857.8
$ws.Range("I126").Value = 2429.6667
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 7289.000100000001
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -4819.000100000001
$ws.Range("N126").Value = -15440
$ws.Range("H132").Value = 3588.5
$ws.Range("I132").Value = 2502.1667
$ws.Range("J132").Value = 4403.25
$ws.Range("K132").Value = 7506.500100000001
$ws.Range("L132").Value = 13209.75
$ws.Range("M132").Value = -4976.500100000001
$ws.Range("N132").Value = -18269.75

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 10007.149
$ws.Range("J131").Value = 10456.169
$ws.Range("L131").Value = 31368.507
$ws.Range("N131").Value = -41448.507
$ws.Range("H132").Value = 1099.1666
$ws.Range("J132").Value = 1500
$ws.Range("L132").Value = 13500
$ws.Range("N132").Value = -18560
$ws.Range("H137").Value = 3469
$ws.Range("I137").Value = 2745.8
$ws.Range("J137").Value = 3985.5715
$ws.Range("K137").Value = 8237.400000000001
$ws.Range("L137").Value = 11956.7145
$ws.Range("M137").Value = -3137.400000000001
$ws.Range("N137").Value = -22156.7145

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H102").Value = 1995.7241
$ws.Range("I102").Value = 2041.4231
$ws.Range("K102").Value = 2041.4231
$ws.Range("M102").Value = -419.4231
$ws.Range("H126").Value = 2461363.2
$ws.Range("I126").Value = 6175495
$ws.Range("J126").Value = 73707.07000000001
$ws.Range("K126").Value = 18526485
$ws.Range("L126").Value = 221121.21
$ws.Range("M126").Value = -18524015
$ws.Range("N126").Value = -226061.21
$ws.Range("H132").Value = 1285356.5
$ws.Range("I132").Value = 1751386.4
$ws.Range("K132").Value = 5254159.199999999
$ws.Range("M132").Value = -5251629.199999999

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 2140.32
$ws.Range("I7").Value = 1456.9131
$ws.Range("K7").Value = 1456.9131
$ws.Range("M7").Value = -1344.9131
$ws.Range("H55").Value = 347.48386
$ws.Range("I55").Value = 282.08
$ws.Range("K55").Value = 282.08
$ws.Range("M55").Value = -109.08
$ws.Range("H93").Value = 886.6
$ws.Range("I93").Value = 886.6
$ws.Range("K93").Value = 886.6
$ws.Range("M93").Value = 361.4
$ws.Range("H126").Value = 2140.32
$ws.Range("I126").Value = 1456.9131
$ws.Range("K126").Value = 4370.7393
$ws.Range("M126").Value = -1900.7393
$ws.Range("H132").Value = 1453.2222
$ws.Range("I132").Value = 1271.2858
$ws.Range("J132").Value = 2090
$ws.Range("K132").Value = 3813.8574
$ws.Range("L132").Value = 6270
$ws.Range("M132").Value = -1283.8574
$ws.Range("N132").Value = -11330
$ws.Range("H136").Value = 3561.6
$ws.Range("I136").Value = 2247.0715
$ws.Range("J136").Value = 5234.636
$ws.Range("K136").Value = 6741.2145
$ws.Range("L136").Value = 15703.908
$ws.Range("M136").Value = -4191.2145
$ws.Range("N136").Value = -20803.908

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value = 130961.836
$ws.Range("I122").Value = 156854.4
$ws.Range("K122").Value = 470563.2
$ws.Range("M122").Value = -468113.2
$ws.Range("H123").Value = 47600
$ws.Range("J123").Value = 47600
$ws.Range("L123").Value = 47600
$ws.Range("N123").Value = -57400
$ws.Range("H132").Value = 2616.3635
$ws.Range("I132").Value = 1256.2858
$ws.Range("J132").Value = 4996.5
$ws.Range("K132").Value = 3768.8574
$ws.Range("L132").Value = 14989.5
$ws.Range("M132").Value = -1238.8574
$ws.Range("N132").Value = -20049.5
$ws.Range("H141").Value = 77899.92999999999
$ws.Range("J141").Value = 77899.92999999999
$ws.Range("L141").Value = 77899.92999999999
$ws.Range("N141").Value = -88259.92999999999
